$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2048929663608563
$ws.Range("C2").Value = 0.5168195718654435
$ws.Range("J2").Value = 0.009174311926605505
$ws.Range("P2").Value = 0.1773700305810398
$ws.Range("S2").Value = 0.09174311926605505

# Row 3
$ws.Range("B3").Value = 0.01104972375690608
$ws.Range("C3").Value = 0.01657458563535912
$ws.Range("J3").Value = 0.02209944751381215
$ws.Range("P3").Value = 0.7458563535911602
$ws.Range("S3").Value = 0.2044198895027624

# Row 4
$ws.Range("J4").Value = 0.04347826086956522
$ws.Range("P4").Value = 0.5217391304347826
$ws.Range("S4").Value = 0.4347826086956522

# Row 6
$ws.Range("B6").Value = 0.05921052631578947
$ws.Range("D6").Value = 0.0131578947368421
$ws.Range("F6").Value = 0.04605263157894737
$ws.Range("J6").Value = 0.2861842105263158
$ws.Range("O6").Value = 0.02302631578947368
$ws.Range("Q6").Value = 0.1513157894736842
$ws.Range("R6").Value = 0.05921052631578947
$ws.Range("S6").Value = 0.3618421052631579

# Row 7
$ws.Range("B7").Value = 0.1644444444444444
$ws.Range("D7").Value = 0.008888888888888889
$ws.Range("F7").Value = 0.06666666666666667
$ws.Range("J7").Value = 0.1422222222222222
$ws.Range("O7").Value = 0.04444444444444445
$ws.Range("Q7").Value = 0.1244444444444444
$ws.Range("R7").Value = 0.04
$ws.Range("S7").Value = 0.4088888888888889

# Row 8
$ws.Range("B8").Value = 0.09696969696969697
$ws.Range("D8").Value = 0.01818181818181818
$ws.Range("F8").Value = 0.05454545454545454
$ws.Range("J8").Value = 0.1272727272727273
$ws.Range("O8").Value = 0.04242424242424243
$ws.Range("Q8").Value = 0.1474747474747475
$ws.Range("R8").Value = 0.06060606060606061
$ws.Range("S8").Value = 0.4525252525252526

# Row 9
$ws.Range("B9").Value = 0.1025641025641026
$ws.Range("D9").Value = 0.02136752136752137
$ws.Range("E9").Value = 0.004273504273504274
$ws.Range("F9").Value = 0.09401709401709402
$ws.Range("J9").Value = 0.1153846153846154
$ws.Range("O9").Value = 0.03846153846153846
$ws.Range("Q9").Value = 0.1324786324786325
$ws.Range("R9").Value = 0.08547008547008547
$ws.Range("S9").Value = 0.405982905982906

# Row 10
$ws.Range("B10").Value = 0.1198120595144871
$ws.Range("D10").Value = 0.02270947533281127
$ws.Range("E10").Value = 0.0007830853563038371
$ws.Range("F10").Value = 0.0908379013312451
$ws.Range("J10").Value = 0.1049334377447142
$ws.Range("O10").Value = 0.03288958496476116
$ws.Range("Q10").Value = 0.1910728269381363
$ws.Range("R10").Value = 0.05168363351605325
$ws.Range("S10").Value = 0.3852779953014879

# Row 11
$ws.Range("G11").Value = 0.1635220125786163
$ws.Range("J11").Value = 0.08490566037735849
$ws.Range("K11").Value = 0.1918238993710692
$ws.Range("L11").Value = 0.5283018867924528
$ws.Range("S11").Value = 0.03144654088050314

# Row 12
$ws.Range("G12").Value = 0.6875
$ws.Range("J12").Value = 0.1428571428571428
$ws.Range("K12").Value = 0.02678571428571428
$ws.Range("L12").Value = 0.0625
$ws.Range("S12").Value = 0.08035714285714286

# Row 13
$ws.Range("F13").Value = 0.01818181818181818
$ws.Range("G13").Value = 0.7090909090909091
$ws.Range("J13").Value = 0.1090909090909091
$ws.Range("S13").Value = 0.1636363636363636

# Row 15
$ws.Range("F15").Value = 0.04642857142857143
$ws.Range("H15").Value = 0.1678571428571428
$ws.Range("I15").Value = 0.06785714285714285
$ws.Range("J15").Value = 0.2821428571428571
$ws.Range("K15").Value = 0.06785714285714285
$ws.Range("M15").Value = 0.007142857142857143
$ws.Range("N15").Value = 0.003571428571428571
$ws.Range("O15").Value = 0.05714285714285714
$ws.Range("S15").Value = 0.3

# Row 16
$ws.Range("F16").Value = 0.03827751196172249
$ws.Range("H16").Value = 0.1961722488038277
$ws.Range("I16").Value = 0.09090909090909091
$ws.Range("J16").Value = 0.354066985645933
$ws.Range("K16").Value = 0.1196172248803828
$ws.Range("M16").Value = 0.02392344497607655
$ws.Range("N16").Value = 0.004784688995215311
$ws.Range("O16").Value = 0.03827751196172249
$ws.Range("S16").Value = 0.1339712918660287

# Row 17
$ws.Range("F17").Value = 0.02147971360381861
$ws.Range("H17").Value = 0.1646778042959427
$ws.Range("I17").Value = 0.09069212410501193
$ws.Range("J17").Value = 0.4295942720763723
$ws.Range("K17").Value = 0.09785202863961814
$ws.Range("M17").Value = 0.02386634844868735
$ws.Range("O17").Value = 0.06921241050119331
$ws.Range("S17").Value = 0.1026252983293556

# Row 18
$ws.Range("F18").Value = 0.02083333333333333
$ws.Range("H18").Value = 0.1597222222222222
$ws.Range("I18").Value = 0.1111111111111111
$ws.Range("J18").Value = 0.3819444444444444
$ws.Range("K18").Value = 0.09027777777777778
$ws.Range("M18").Value = 0.02083333333333333
$ws.Range("O18").Value = 0.08333333333333333
$ws.Range("S18").Value = 0.1319444444444444

# Row 19
$ws.Range("F19").Value = 0.04016337644656229
$ws.Range("H19").Value = 0.20285908781484
$ws.Range("I19").Value = 0.0878148400272294
$ws.Range("J19").Value = 0.3451327433628318
$ws.Range("K19").Value = 0.09530292716133425
$ws.Range("M19").Value = 0.01837985023825732
$ws.Range("N19").Value = 0.002042205582028591
$ws.Range("O19").Value = 0.07147719537100068
$ws.Range("S19").Value = 0.1368277739959156
